$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full measures 2000-01-31")

# Swap the "Pearson" / "Spearman" measure rows (row 4 <-> row 5, columns P:R).
# Before: row 4 = Pearson  / PEARSON(...) / 0.618986
#         row 5 = Spearman / CORREL(...)  / 0.639752
# After:  row 4 = Spearman / CORREL(...)  / 0.639752
#         row 5 = Pearson  / PEARSON(...) / 0.618986

$ws.Range("P4").Value = "Spearman"
$ws.Range("Q4").Formula = "=CORREL(J4:J25,K4:K25)"

$ws.Range("P5").Value = "Pearson"
$ws.Range("Q5").Formula = "=PEARSON(H4:H25,I4:I25)"

# R4/R5 hold their numbers as text (shared strings), so force text storage
# before writing, then drop the formatting again to match the original
# (unstyled) cells.
$ws.Range("R4:R5").NumberFormat = "@"
$ws.Range("R4").Value = "0.639752"
$ws.Range("R5").Value = "0.618986"
$ws.Range("R4:R5").ClearFormats()

# P11 used to be an empty-but-formatted cell - remove it (content + format) entirely.
$ws.Range("P11").Clear()

# Update the saved selection.
$ws.Range("K13").Select()
